# Update the "Results" table on slide 12 ("Tabella 10"):
#   - Logistic Regression / Numerical / with standardization: 0.977 -> 0.9772
#   - Logistic Regression / Numerical / w/out standardize:    (empty) -> 0.9773
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$shp = $s.Shapes.Item("Tabella 10")
$tbl = $shp.Table

# Row 3 = "Logistic Regression", Column 2 = "Numerical / with standardization"
$cellWith = $tbl.Cell(3, 2)
$cellWith.Shape.TextFrame.TextRange.Text = "0.9772"

# Row 3 = "Logistic Regression", Column 3 = "Numerical / w/out standardize" (was empty)
$cellWithout = $tbl.Cell(3, 3)
$trWithout = $cellWithout.Shape.TextFrame.TextRange
$trWithout.Text = "0.9773"
$trWithout.Font.Name = "Avenir Next LT Pro"
